# "added script for add loan type screen"
# Extends the Institute sheet (sheet10.xml) with 5 new trailing columns
# (BL:BP -> HolidayType, LoanType, DraftNeeded, LoanTypeCreditLimit,
# LoanTypeCashLimit) and fills them in for the 11 existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Institute")

# --- copy formatting for the new columns from the last existing column ---
# Row 1 (header) uses the same style as BK1; rows 2-12 (data) use the same
# style as BK2 (identical for every data row in this sheet).
$ws.Range("BK1").Copy() | Out-Null
$ws.Range("BL1:BP1").PasteSpecial(-4122) | Out-Null

$ws.Range("BK2").Copy() | Out-Null
$ws.Range("BL2:BP12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- header row ---
$ws.Range("BL1").Value = "HolidayType"
$ws.Range("BM1").Value = "LoanType"
$ws.Range("BN1").Value = "DraftNeeded"
$ws.Range("BO1").Value = "LoanTypeCreditLimit"
$ws.Range("BP1").Value = "LoanTypeCashLimit"

# --- data rows ---
# Row 3 has a different LoanType value than the rest of the rows.
$rows = 2..12
foreach ($r in $rows) {
    $ws.Range("BL$r").Value = "Holiday [H]"
    if ($r -eq 3) {
        $ws.Range("BM$r").Value = "Retail Transaction to Loan [LOANPUR]"
    } else {
        $ws.Range("BM$r").Value = "Loan below credit limit [LOANCR]"
    }
    $ws.Range("BN$r").Value = "Check"
    $ws.Range("BO$r").Value = "Check"
    $ws.Range("BP$r").Value = "Check"
}

# --- column width for the new last column (BP) ---
$ws.Columns.Item(68).ColumnWidth = 17.5

# --- view state: scroll position + active selection ---
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 55
$win.ScrollRow = 1
$ws.Range("BO4").Select()
